$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row before the current row 194; this pushes the
# existing rows 194-302 down to 195-303 and extends the sheet dimension
# from A1:R302 to A1:R303.
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new weekly price record.
$ws.Range("A194").Value = 7
$ws.Range("B194").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C194").Value = "Ñuble"
$ws.Range("D194").Value = 44488
$ws.Range("E194").Value = 16
$ws.Range("F194").Value = 100112020
$ws.Range("G194").Value = "Tomate"
$ws.Range("H194").Value = "Larga vida"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 360
$ws.Range("K194").Value = 7500
$ws.Range("L194").Value = 8000
$ws.Range("M194").Value = 7750
$ws.Range("N194").Value = "$/caja 10 kilos"
$ws.Range("O194").Value = "Región de Arica y Parinacota"
$ws.Range("P194").Value = 775
$ws.Range("Q194").Value = 10
$ws.Range("R194").Value = "Hortaliza"
